# Append 16 new "FAIL" rows (949-964) to Sheet1, mirroring the next daily-update
# upload: new Circle / Site Id / Status Date / Site Status / Fail KPI / Test Case /
# Remarks entries appended right after the previous last row (948).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Template row whose formatting (borders/fonts/wrap, incl. the red "FAIL" style
# in column D) every new row below reuses.
$templateRow = 948

# Row 949: MPBP6157_BHO_P40
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A949:G949").PasteSpecial(-4122)
$ws.Cells.Item(949,1).Value = "MP"
$ws.Cells.Item(949,2).Value = "MPBP6157_BHO_P40"
$ws.Cells.Item(949,3).Value = "22-Dec-2025 9:54 AM"
$ws.Cells.Item(949,4).Value = "FAIL"
$ws.Cells.Item(949,5).Value = "1. Peak PUSCH UL Throughput"
$ws.Cells.Item(949,6).Value = "1. Static UL"
$ws.Cells.Item(949,7).Value = "1. Peak PUSCH DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the maximum value of PUSCH Throughput in the NR tab."
$ws.Rows.Item(949).RowHeight = 16

# Row 950: TOND49_CHN_P41
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A950:G950").PasteSpecial(-4122)
$ws.Cells.Item(950,1).Value = "CH"
$ws.Cells.Item(950,2).Value = "TOND49_CHN_P41"
$ws.Cells.Item(950,3).Value = "22-Dec-2025 12:55 AM"
$ws.Cells.Item(950,4).Value = "FAIL"
$ws.Cells.Item(950,5).Value = "1. MO Call (pass/fail)`n2. MT Call (pass/fail)`n3. CSFB Call (pass/fail)`n4. Ping/Round trip time(ms)"
$ws.Cells.Item(950,6).Value = "1. Static VoLTE MO`n2. Static VoLTE MT`n3. Static CSFB MO`n4. Static Ping"
$ws.Cells.Item(950,7).Value = "1. VoLTE Long Call MO – As per Bharti acceptance criteria, a minimum of 3 successful call setups are required without any blocked call. Kindly perform at least 3 successful MO calls.`n2. VoLTE Long Call MT – As per Bharti acceptance criteria, a minimum of 3 successful call setups are required without any blocked call. Kindly perform at least 3 successful MT calls.`n3. CSFB MO – As per Bharti acceptance criteria, a minimum of 3 successful call setups are required without any blocked call. Kindly perform at least 3 successful MO  calls.`n4. Ping is not meeting the acceptance criteria. The average ping value across all logfiles should be less than 50 ms. Kindly exclude the logfile where the average value exceeds 50 ms and redo the test.”"
$ws.Rows.Item(950).RowHeight = 64

# Row 951: TIR644_TIR_P40
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A951:G951").PasteSpecial(-4122)
$ws.Cells.Item(951,1).Value = "TN"
$ws.Cells.Item(951,2).Value = "TIR644_TIR_P40"
$ws.Cells.Item(951,3).Value = "21-Dec-2025 10:16 PM"
$ws.Cells.Item(951,4).Value = "FAIL"
$ws.Cells.Item(951,5).Value = "1. Video Streaming  (ms)"
$ws.Cells.Item(951,6).Value = "1. Static Yotube Streaming"
$ws.Cells.Item(951,7).Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$ws.Rows.Item(951).RowHeight = 24

# Row 952: MDU712_MDR_P40
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A952:G952").PasteSpecial(-4122)
$ws.Cells.Item(952,1).Value = "TN"
$ws.Cells.Item(952,2).Value = "MDU712_MDR_P40"
$ws.Cells.Item(952,3).Value = "21-Dec-2025 8:45 PM"
$ws.Cells.Item(952,4).Value = "FAIL"
$ws.Cells.Item(952,5).Value = "1. Video Streaming  (ms)"
$ws.Cells.Item(952,6).Value = "1. Static Yotube Streaming"
$ws.Cells.Item(952,7).Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$ws.Rows.Item(952).RowHeight = 24

# Row 953: BHPAT-1287_PAT_P41
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A953:G953").PasteSpecial(-4122)
$ws.Cells.Item(953,1).Value = "BH"
$ws.Cells.Item(953,2).Value = "BHPAT-1287_PAT_P41"
$ws.Cells.Item(953,3).Value = "21-Dec-2025 7:42 PM"
$ws.Cells.Item(953,4).Value = "FAIL"
$ws.Cells.Item(953,5).Value = "1. SCG addition after VoLTE call released`n2. SgNB Addition time (ms)`n3. Video Streaming"
$ws.Cells.Item(953,6).Value = "1. Static VoLTE MO`n2. Static ATDT`n3. Static Yotube Streaming"
$ws.Cells.Item(953,7).Value = "1. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition.`n2. Exclue ATDT Logfile and Create New Sgnb Addition Time Is Very High. It Should Be <150 Ms. To Achieve This, Perform Static Test In Main Lobe And Keep Test Files Downloading In Background. Also, Ensure 4G Serving Cell Belongs To The Same Site. Exclude The Existing Logfile First`n3. While performing the YouTube test for both sectors, please ensure that the video is playing successfully in the script before saving the log file."
$ws.Rows.Item(953).RowHeight = 72

# Row 954: PUN6448_PNE_P40
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A954:G954").PasteSpecial(-4122)
$ws.Cells.Item(954,1).Value = "MH"
$ws.Cells.Item(954,2).Value = "PUN6448_PNE_P40"
$ws.Cells.Item(954,3).Value = "21-Dec-2025 4:23 PM"
$ws.Cells.Item(954,4).Value = "FAIL"
$ws.Cells.Item(954,5).Value = "1. Video Streaming  (ms)"
$ws.Cells.Item(954,6).Value = "1. Static Yotube Streaming"
$ws.Cells.Item(954,7).Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$ws.Rows.Item(954).RowHeight = 24

# Row 955: KUR5146_SLP_P40
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A955:G955").PasteSpecial(-4122)
$ws.Cells.Item(955,1).Value = "MH"
$ws.Cells.Item(955,2).Value = "KUR5146_SLP_P40"
$ws.Cells.Item(955,3).Value = "21-Dec-2025 10:03 AM"
$ws.Cells.Item(955,4).Value = "FAIL"
$ws.Cells.Item(955,5).Value = "1. Video Streaming  (ms)"
$ws.Cells.Item(955,6).Value = "1. Static Yotube Streaming"
$ws.Cells.Item(955,7).Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$ws.Rows.Item(955).RowHeight = 24

# Row 956: KDTP15_PKD_P40
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A956:G956").PasteSpecial(-4122)
$ws.Cells.Item(956,1).Value = "KL"
$ws.Cells.Item(956,2).Value = "KDTP15_PKD_P40"
$ws.Cells.Item(956,3).Value = "20-Dec-2025 9:39 PM"
$ws.Cells.Item(956,4).Value = "FAIL"
$ws.Cells.Item(956,5).Value = "1. SCG addition after VoLTE call released`n2. Peak PDSCH DL Throughput"
$ws.Cells.Item(956,6).Value = "1. Static VoLTE MO`n2. Static DL"
$ws.Cells.Item(956,7).Value = "1. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition.`n2. Peak PDSCH DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the maximum value of PDSCH Throughput in the NR tab."
$ws.Rows.Item(956).RowHeight = 48

# Row 957: AJJP02_JPR_P40
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A957:G957").PasteSpecial(-4122)
$ws.Cells.Item(957,1).Value = "RJ"
$ws.Cells.Item(957,2).Value = "AJJP02_JPR_P40"
$ws.Cells.Item(957,3).Value = "20-Dec-2025 8:22 PM"
$ws.Cells.Item(957,4).Value = "FAIL"
$ws.Cells.Item(957,5).Value = "1. CSFB Call (pass/fail)"
$ws.Cells.Item(957,6).Value = "1. Static CSFB MO"
$ws.Cells.Item(957,7).Value = "1. CSFB MO – As per Bharti acceptance criteria, a minimum of 3 successful call setups are required without any blocked call. Kindly perform at least 3 successful MO  calls."
$ws.Rows.Item(957).RowHeight = 16

# Row 958: BHFOR-29_PAT_P40
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A958:G958").PasteSpecial(-4122)
$ws.Cells.Item(958,1).Value = "BH"
$ws.Cells.Item(958,2).Value = "BHFOR-29_PAT_P40"
$ws.Cells.Item(958,3).Value = "20-Dec-2025 7:13 PM"
$ws.Cells.Item(958,4).Value = "FAIL"
$ws.Cells.Item(958,5).Value = "1. SCG addition after VoLTE call released"
$ws.Cells.Item(958,6).Value = "1. Static VoLTE MO"
$ws.Cells.Item(958,7).Value = "1. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition."
$ws.Rows.Item(958).RowHeight = 32

# Row 959: KOZ568_KOZ_P40
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A959:G959").PasteSpecial(-4122)
$ws.Cells.Item(959,1).Value = "KL"
$ws.Cells.Item(959,2).Value = "KOZ568_KOZ_P40"
$ws.Cells.Item(959,3).Value = "20-Dec-2025 5:04 PM"
$ws.Cells.Item(959,4).Value = "FAIL"
$ws.Cells.Item(959,5).Value = "1. Peak PDSCH DL Throughput`n2. Video Streaming  (ms)"
$ws.Cells.Item(959,6).Value = "1. Static DL`n2. Static Yotube Streaming"
$ws.Cells.Item(959,7).Value = "1. Peak PDSCH DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the maximum value of PDSCH Throughput in the NR tab.`n2. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$ws.Rows.Item(959).RowHeight = 40

# Row 960: BHMWJ-02_1_PAT_P40
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A960:G960").PasteSpecial(-4122)
$ws.Cells.Item(960,1).Value = "BH"
$ws.Cells.Item(960,2).Value = "BHMWJ-02_1_PAT_P40"
$ws.Cells.Item(960,3).Value = "20-Dec-2025 2:33 PM"
$ws.Cells.Item(960,4).Value = "FAIL"
$ws.Cells.Item(960,5).Value = "1. MT Call (pass/fail)`n2. SCG addition after VoLTE call released`n3. Peak PDCP DL Throughput`n4. Average PDCP DL Throughput`n5. Peak PDSCH DL Throughput`n6. Median PDCP DL​ Throughput`n7. Downlink Peak MCS - 5G`n8. Peak PDCP UL Throughput`n9. Avg PDCP UL Throughput`n10. Peak PUSCH UL Throughput`n11. Median PDCP UL Throughput`n12. SgNB Addition time (ms)`n13. UE Steering (Idle) : Non anchor/anchor to preferred anchor"
$ws.Cells.Item(960,6).Value = "1. Static VoLTE MT`n2. Static VoLTE MO`n3. Static DL`n4. Static DL`n5. Static DL`n6. Mobility DL`n7. Static DL`n8. Static UL`n9. Static UL`n10. Static UL`n11. Mobility UL`n12. Static ATDT`n13. Static Idle"
$ws.Cells.Item(960,7).Value = "1. VoLTE Long Call MT – As per Bharti acceptance criteria, a minimum of 3 successful call setups are required without any blocked call. Kindly perform at least 3 successful MT calls.`n2. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition.`n3. Peak PDCP DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the reporting of MR-DC DRB PDCP DL Throughput in the NR tab.`n4. Average PDCP DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the reporting of MR-DC DRB PDCP DL Throughput in the NR tab.`n5. Peak PDSCH DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the maximum value of PDSCH Throughput in the NR tab.`n6. The Median PDCP DL Throughput is reported as 0. Kindly add or exclude a logfile in the DL drive so the median value can update. It is recommended to add a new logfile and collect maximum throughput samples in a good coverage area.`n7. Peak MCS is not meeting the acceptance criteria. Kindly redo the test and verify that the value meets the required threshold. To achieve the desired MCS, perform the test in the main lobe of the cell within a good coverage area.`n8. Peak PDCP UL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the reporting of MR-DC DRB PDCP UL Throughput in the NR tab.`n9. Average PDCP UL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the reporting of MR-DC DRB PDCP UL Throughput in the NR tab.`n10. Peak PUSCH DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the maximum value of PUSCH Throughput in the NR tab.`n11. The Median PDCP UL Throughput is reported as 0. Kindly add or exclude a logfile in the DL drive so the median value can update. It is recommended to add a new logfile and collect maximum throughput samples in a good coverage area.`n12. Exclue ATDT Logfile and Create New Sgnb Addition Time Is Very High. It Should Be <150 Ms. To Achieve This, Perform Static Test In Main Lobe And Keep Test Files Downloading In Background. Also, Ensure 4G Serving Cell Belongs To The Same Site. Exclude The Existing Logfile First`n13. For sites with NOKIA OEM, validate using Drive Idle, and for other OEMs, validate using Static Idle. In both Drive and Static Idle, the UE should latch from NR to LTE and from LTE to NR. In LTE, the UE should latch on the band that corresponds to the configured anchor layer."
$ws.Rows.Item(960).RowHeight = 280

# Row 961: CHEL28_1_MLP_P40
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A961:G961").PasteSpecial(-4122)
$ws.Cells.Item(961,1).Value = "KL"
$ws.Cells.Item(961,2).Value = "CHEL28_1_MLP_P40"
$ws.Cells.Item(961,3).Value = "20-Dec-2025 2:48 PM"
$ws.Cells.Item(961,4).Value = "FAIL"
$ws.Cells.Item(961,5).Value = "1. SCG addition after VoLTE call released`n2. Downlink Peak MCS - 5G`n3. SgNB Addition time (ms)`n4. Video Streaming  (ms)"
$ws.Cells.Item(961,6).Value = "1. Static VoLTE MO`n2. Static DL`n3. Static ATDT`n4. Static Yotube Streaming"
$ws.Cells.Item(961,7).Value = "1. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition.`n2. Peak MCS is not meeting the acceptance criteria. Kindly redo the test and verify that the value meets the required threshold. To achieve the desired MCS, perform the test in the main lobe of the cell within a good coverage area.`n3. Exclue ATDT Logfile and Create New Sgnb Addition Time Is Very High. It Should Be <150 Ms. To Achieve This, Perform Static Test In Main Lobe And Keep Test Files Downloading In Background. Also, Ensure 4G Serving Cell Belongs To The Same Site. Exclude The Existing Logfile First`n4. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$ws.Rows.Item(961).RowHeight = 96

# Row 962: MDU747_MDR_P40
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A962:G962").PasteSpecial(-4122)
$ws.Cells.Item(962,1).Value = "TN"
$ws.Cells.Item(962,2).Value = "MDU747_MDR_P40"
$ws.Cells.Item(962,3).Value = "20-Dec-2025 12:13 PM"
$ws.Cells.Item(962,4).Value = "FAIL"
$ws.Cells.Item(962,5).Value = "1. Video Streaming  (ms)"
$ws.Cells.Item(962,6).Value = "1. Static Yotube Streaming"
$ws.Cells.Item(962,7).Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$ws.Rows.Item(962).RowHeight = 24

# Row 963: BHMJQ-02_1_PAT_P40
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A963:G963").PasteSpecial(-4122)
$ws.Cells.Item(963,1).Value = "BH"
$ws.Cells.Item(963,2).Value = "BHMJQ-02_1_PAT_P40"
$ws.Cells.Item(963,3).Value = "20-Dec-2025 11:35 AM"
$ws.Cells.Item(963,4).Value = "FAIL"
$ws.Cells.Item(963,5).Value = "1. PCI`n2. SCG addition after VoLTE call released"
$ws.Cells.Item(963,6).Value = "1. Static All`n2. Static VoLTE MO"
$ws.Cells.Item(963,7).Value = "1. The PCI uploaded in the site database is not aligned with the actual on-site servings. Kindly verify the configured PCI and update the database accordingly to avoid inconsistency in reporting.`n2. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition."
$ws.Rows.Item(963).RowHeight = 48

# Row 964: PUN6450_PNE_P40
$ws.Range("A" + $templateRow + ":G" + $templateRow).Copy()
$ws.Range("A964:G964").PasteSpecial(-4122)
$ws.Cells.Item(964,1).Value = "MH"
$ws.Cells.Item(964,2).Value = "PUN6450_PNE_P40"
$ws.Cells.Item(964,3).Value = "20-Dec-2025 10:53 AM"
$ws.Cells.Item(964,4).Value = "FAIL"
$ws.Cells.Item(964,5).Value = "1. Video Streaming  (ms)"
$ws.Cells.Item(964,6).Value = "1. Static Yotube Streaming"
$ws.Cells.Item(964,7).Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$ws.Rows.Item(964).RowHeight = 24

# Restore the selection state left by the author after this edit session
# (row-height/AutoFit UI feedback is not modeled by this headless host, so the
# scroll position (topLeftCell) cannot be replicated here).
$ws.Range("E973").Select()
